# Weekly fruit/vegetable price update: insert a new daily record as row 46
# ("Pepino dulce" / Vega Monumental Concepción), pushing the existing
# rows 46-75 down to 47-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 46 (copies formatting, incl. the
# date style on column D, from the row above - matching native Excel
# "Insert Row" behaviour).
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new market record.
$ws.Cells.Item(46, 1).Value  = 11
$ws.Cells.Item(46, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value  = "Bíobío"
$ws.Cells.Item(46, 4).Value  = 45086
$ws.Cells.Item(46, 5).Value  = 8
$ws.Cells.Item(46, 6).Value  = 100112043
$ws.Cells.Item(46, 7).Value  = "Pepino dulce"
$ws.Cells.Item(46, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(46, 9).Value  = "Primera"
$ws.Cells.Item(46, 10).Value = 200
$ws.Cells.Item(46, 11).Value = 13000
$ws.Cells.Item(46, 12).Value = 14000
$ws.Cells.Item(46, 13).Value = 13500
$ws.Cells.Item(46, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(46, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value = 750
$ws.Cells.Item(46, 17).Value = 18
$ws.Cells.Item(46, 18).Value = "Hortaliza"
